$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 109.333336
$ws.Range("I33").Value = 69.14286
$ws.Range("J33").Value = 390.66666
$ws.Range("K33").Value = 69.14286
$ws.Range("L33").Value = 390.66666
$ws.Range("M33").Value = 159.85714
$ws.Range("N33").Value = -848.66666
# Row 62
$ws.Range("H62").Value = 4273.9473
$ws.Range("I62").Value = 2518.6365
$ws.Range("J62").Value = 6687.5
$ws.Range("K62").Value = 2518.6365
$ws.Range("L62").Value = 6687.5
$ws.Range("M62").Value = -1894.6365
$ws.Range("N62").Value = -7935.5
# Row 65
$ws.Range("H65").Value = 4273.9473
$ws.Range("I65").Value = 2518.6365
$ws.Range("J65").Value = 6687.5
$ws.Range("K65").Value = 12593.1825
$ws.Range("L65").Value = 33437.5
$ws.Range("M65").Value = -9473.182500000001
$ws.Range("N65").Value = -39677.5
# Row 100
$ws.Range("H100").Value = 1488.8518
$ws.Range("I100").Value = 1341.1052
$ws.Range("K100").Value = 1341.1052
$ws.Range("M100").Value = -800.1052
# Row 113
$ws.Range("H113").Value = 2828
$ws.Range("I113").Value = 1868.3334
$ws.Range("J113").Value = 3147.889
$ws.Range("K113").Value = 1868.3334
$ws.Range("L113").Value = 3147.889
$ws.Range("M113").Value = 1385.6666
$ws.Range("N113").Value = -9655.888999999999
# Row 130
$ws.Range("H130").Value = 45179.547
$ws.Range("J130").Value = 45179.547
$ws.Range("L130").Value = 45179.547
$ws.Range("N130").Value = -55219.547
# Row 138
$ws.Range("H138").Value = 5197.0977
$ws.Range("I138").Value = 1595.5714
$ws.Range("J138").Value = 8978.700000000001
$ws.Range("K138").Value = 4786.7142
$ws.Range("L138").Value = 26936.1
$ws.Range("M138").Value = 353.2857999999997
$ws.Range("N138").Value = -37216.10000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2642.3076
$ws.Range("I2").Value = 3415.8572
$ws.Range("J2").Value = 1739.8334
$ws.Range("K2").Value = 3415.8572
$ws.Range("L2").Value = 1739.8334
$ws.Range("M2").Value = -3302.8572
$ws.Range("N2").Value = -1965.8334
# Row 110
$ws.Range("H110").Value = 1661
$ws.Range("I110").Value = 1492.4445
$ws.Range("K110").Value = 1492.4445
$ws.Range("M110").Value = 552.5554999999999
# Row 116
$ws.Range("H116").Value = 2642.3076
$ws.Range("I116").Value = 3415.8572
$ws.Range("J116").Value = 1739.8334
$ws.Range("K116").Value = 3415.8572
$ws.Range("L116").Value = 1739.8334
$ws.Range("M116").Value = -1121.8572
$ws.Range("N116").Value = -6327.8334
# Row 132
$ws.Range("H132").Value = 2064.6155
$ws.Range("I132").Value = 1875.5714
$ws.Range("J132").Value = 2545.818
$ws.Range("K132").Value = 5626.7142
$ws.Range("L132").Value = 7637.454000000001
$ws.Range("M132").Value = -3096.7142
$ws.Range("N132").Value = -12697.454

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2642.3076
$ws.Range("I3").Value = 3415.8572
$ws.Range("J3").Value = 1739.8334
$ws.Range("K3").Value = 3415.8572
$ws.Range("L3").Value = 1739.8334
$ws.Range("M3").Value = -3301.8572
$ws.Range("N3").Value = -1967.8334
# Row 22
$ws.Range("H22").Value = 263.66666
# Row 104
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
# Row 107
$ws.Range("H107").Value = 1180.55
$ws.Range("I107").Value = 969.4375
$ws.Range("J107").Value = 2025
$ws.Range("K107").Value = 969.4375
$ws.Range("L107").Value = 2025
$ws.Range("M107").Value = 950.5625
$ws.Range("N107").Value = -5865
# Row 134
$ws.Range("H134").Value = 2687.15
$ws.Range("I134").Value = 2766.2727
$ws.Range("J134").Value = 2314.1428
$ws.Range("K134").Value = 8298.8181
$ws.Range("L134").Value = 6942.428400000001
$ws.Range("M134").Value = -5763.8181
$ws.Range("N134").Value = -12012.4284

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3682.6
$ws.Range("I16").Value = 2966.6667
$ws.Range("J16").Value = 3989.4285
$ws.Range("K16").Value = 2966.6667
$ws.Range("L16").Value = 3989.4285
$ws.Range("M16").Value = -2679.6667
$ws.Range("N16").Value = -4563.4285
# Row 107
$ws.Range("H107").Value = 836.4286
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 839.2308
$ws.Range("K107").Value = 800
$ws.Range("L107").Value = 839.2308
$ws.Range("M107").Value = 1120
$ws.Range("N107").Value = -4679.2308
# Row 109
$ws.Range("H109").Value = 46642.5
$ws.Range("J109").Value = 46642.5
$ws.Range("L109").Value = 46642.5
$ws.Range("N109").Value = -48722.5
# Row 113
$ws.Range("H113").Value = 3682.6
$ws.Range("I113").Value = 2966.6667
$ws.Range("J113").Value = 3989.4285
$ws.Range("K113").Value = 2966.6667
$ws.Range("L113").Value = 3989.4285
$ws.Range("M113").Value = -796.6667000000002
$ws.Range("N113").Value = -8329.4285

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 7122.45
$ws.Range("I87").Value = 885.1667
$ws.Range("J87").Value = 9795.571
$ws.Range("K87").Value = 2655.5001
$ws.Range("L87").Value = 29386.713
$ws.Range("M87").Value = -1407.5001
$ws.Range("N87").Value = -31882.713
# Row 90
$ws.Range("H90").Value = 7122.45
$ws.Range("I90").Value = 885.1667
$ws.Range("J90").Value = 9795.571
$ws.Range("K90").Value = 7966.5003
$ws.Range("L90").Value = 88160.139
$ws.Range("M90").Value = -1726.5003
$ws.Range("N90").Value = -100640.139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 419.375
$ws.Range("I107").Value = 150.28572
$ws.Range("J107").Value = 628.6667
$ws.Range("K107").Value = 150.28572
$ws.Range("L107").Value = 628.6667
$ws.Range("M107").Value = 1769.71428
$ws.Range("N107").Value = -4468.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 767.5
$ws.Range("I22").Value = 499.66666
$ws.Range("J22").Value = 1035.3334
$ws.Range("K22").Value = 499.66666
$ws.Range("L22").Value = 1035.3334
$ws.Range("M22").Value = -204.66666
$ws.Range("N22").Value = -1625.3334
# Row 27
$ws.Range("H27").Value = 767.5
$ws.Range("I27").Value = 499.66666
$ws.Range("J27").Value = 1035.3334
$ws.Range("K27").Value = 499.66666
$ws.Range("L27").Value = 1035.3334
$ws.Range("M27").Value = -392.66666
$ws.Range("N27").Value = -1249.3334
# Row 61
$ws.Range("H61").Value = 1276988.5
$ws.Range("I61").Value = 35150.5
$ws.Range("J61").Value = 5002502.5
$ws.Range("K61").Value = 35150.5
$ws.Range("L61").Value = 5002502.5
$ws.Range("M61").Value = -34948.5
$ws.Range("N61").Value = -5002906.5
# Row 94
$ws.Range("H94").Value = 21875
$ws.Range("J94").Value = 21875
$ws.Range("L94").Value = 21875
$ws.Range("N94").Value = -23227
# Row 113
$ws.Range("H113").Value = 1276988.5
$ws.Range("I113").Value = 35150.5
$ws.Range("J113").Value = 5002502.5
$ws.Range("K113").Value = 35150.5
$ws.Range("L113").Value = 5002502.5
$ws.Range("M113").Value = -32980.5
$ws.Range("N113").Value = -5006842.5
# Row 132
$ws.Range("H132").Value = 5616.5
$ws.Range("I132").Value = 5680
$ws.Range("J132").Value = 5571.143
$ws.Range("K132").Value = 17040
$ws.Range("L132").Value = 16713.429
$ws.Range("M132").Value = -14510
$ws.Range("N132").Value = -21773.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2917.2856
$ws.Range("I107").Value = 547
$ws.Range("J107").Value = 6077.6665
$ws.Range("K107").Value = 1641
$ws.Range("L107").Value = 18232.9995
$ws.Range("M107").Value = 279
$ws.Range("N107").Value = -22072.9995
# Row 113
$ws.Range("H113").Value = 655.6857
$ws.Range("I113").Value = 377.77777
$ws.Range("J113").Value = 949.94116
$ws.Range("K113").Value = 1133.33331
$ws.Range("L113").Value = 2849.82348
$ws.Range("M113").Value = 1036.66669
$ws.Range("N113").Value = -7189.82348
# Row 135
$ws.Range("H135").Value = 47333.332
$ws.Range("J135").Value = 47333.332
$ws.Range("L135").Value = 47333.332
$ws.Range("N135").Value = -57473.332
